# Fruta / hortaliza, semanal
#
# A new week of data is inserted into rows 6-7 (dated 44559, with updated
# prices). The previously-latest week (dated 44195) that used to occupy
# rows 6-7 is preserved, unchanged, by pushing it down to rows 8-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Range("D6").NumberFormat

# Row 8: duplicate of the old row 6 (Primera, week of 44195).
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44195
$ws.Range("D8").NumberFormat = $dateFormat
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101004
$ws.Range("J8").Value = "Frambuesa"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 3000
$ws.Range("O8").Value = 3500
$ws.Range("P8").Value = 3250
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Región de Ñuble"
$ws.Range("S8").Value = 1625
$ws.Range("T8").Value = 2

# Row 9: duplicate of the old row 7 (Segunda, week of 44195).
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44195
$ws.Range("D9").NumberFormat = $dateFormat
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101004
$ws.Range("J9").Value = "Frambuesa"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 2500
$ws.Range("O9").Value = 2500
$ws.Range("P9").Value = 2500
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Región de Ñuble"
$ws.Range("S9").Value = 1250
$ws.Range("T9").Value = 2

# Row 6: overwrite with the new week's Primera values.
$ws.Range("D6").Value = 44559
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6500
$ws.Range("S6").Value = 3250

# Row 7: overwrite with the new week's Segunda values.
$ws.Range("D7").Value = 44559
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 5000
$ws.Range("P7").Value = 5000
$ws.Range("S7").Value = 2500
